# Apply the "modified and test update" edit to the employees worksheet.
# - Row 3 (employeeid 1001) gets new deparment/image/role/address/qualification.
# - Row 4 (employeeid 1002) gets a new name/email/image (Kajol).
# - Four brand-new rows (5-8 / employeeid 1003-1006) are appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 3: update deparment/image/role/address/qualification ----
$ws.Range("D3").Value = "fynd Accedamy"
$ws.Range("E3").Value = "http://res.cloudinary.com/db5vn6bj5/image/upload/v1635076447/durz7fokpsouywmymjib.jpg"
$ws.Range("F3").Value = "Experiance"
$ws.Range("G3").Value = "Navi mumbai"
$ws.Range("H3").Value = "BSC"

# ---- Row 4: update Name/email/image ----
$ws.Range("B4").Value = "Kajol"
$ws.Range("C4").Value = "shreyajaiswal1996vns@gmail.com"
$ws.Range("E4").Value = "http://res.cloudinary.com/db5vn6bj5/image/upload/v1635076339/aqod8orn4pcytzivegoe.jpg"

# ---- Row 5 (new): employeeid 1003, Deepanshu Gupta ----
$ws.Range("A5").Value = 1003
$ws.Range("B5").Value = "Deepanshu Gupta"
$ws.Range("C5").Value = "tpo@abesit.edu.com"
$ws.Range("D5").Value = "jio mart"
$ws.Range("E5").Value = "http://res.cloudinary.com/db5vn6bj5/image/upload/v1635076503/bgzcctj8mrxpljmd9a09.jpg"
$ws.Range("F5").Value = "Fresher"
$ws.Range("G5").Value = "Rampur Mumbai "
$ws.Range("H5").Value = "B.tech"
$ws.Range("I5").Value = "full time"
$ws.Range("J5").Value = 44493.508656226855
$ws.Range("K5").Value = 123456789

# ---- Row 6 (new): employeeid 1004, Khushabu Belsare ----
$ws.Range("A6").Value = 1004
$ws.Range("B6").Value = "Khushabu Belsare"
$ws.Range("C6").Value = "khushabubelsare68@gmail.com"
$ws.Range("D6").Value = "jio mart"
$ws.Range("E6").Value = "http://res.cloudinary.com/db5vn6bj5/image/upload/v1635076475/vrzazbmaq5ozrfdnzdsr.jpg"
$ws.Range("F6").Value = "Fresher"
$ws.Range("G6").Value = "kolkata"
$ws.Range("H6").Value = "B.tech"
$ws.Range("I6").Value = "full time"
$ws.Range("J6").Value = 44493.508656226855
$ws.Range("K6").Value = 123456789

# ---- Row 7 (new): employeeid 1005, Rahul Tambe ----
$ws.Range("A7").Value = 1005
$ws.Range("B7").Value = "Rahul Tambe"
$ws.Range("C7").Value = "iamrahultambe@gmail.com"
$ws.Range("D7").Value = "jio mart"
$ws.Range("E7").Value = "http://res.cloudinary.com/db5vn6bj5/image/upload/v1635076407/hbdxeoukkn7jl8bavygu.jpg"
$ws.Range("F7").Value = "Experiance"
$ws.Range("G7").Value = "kolkata"
$ws.Range("H7").Value = "B.tech"
$ws.Range("I7").Value = "full time"
$ws.Range("J7").Value = 44493.508656226855
$ws.Range("K7").Value = 123456789

# ---- Row 8 (new): employeeid 1006, Neha Jaiswal ----
$ws.Range("A8").Value = 1006
$ws.Range("B8").Value = "Neha Jaiswal"
$ws.Range("C8").Value = "nehajaiswal694@gmail.com"
$ws.Range("D8").Value = "jio mart"
$ws.Range("E8").Value = "http://res.cloudinary.com/db5vn6bj5/image/upload/v1635076300/hf4jdrb5p2zh1sgnznuf.jpg"
$ws.Range("F8").Value = "Experiance"
$ws.Range("G8").Value = "kolkata"
$ws.Range("H8").Value = "MCA"
$ws.Range("I8").Value = "full time"
$ws.Range("J8").Value = 44493.508656226855
$ws.Range("K8").Value = 123456789

# Copy the date formatting (style index) from the existing J4 cell onto the
# newly added date cells so they share the same number format as the rest
# of the dateOfJoin column instead of getting a brand new style entry.
$ws.Range("J4").Copy()
$ws.Range("J5:J8").PasteSpecial(-4122)
$excel.CutCopyMode = 0
